$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.757.50"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.140.90"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.61"
$ws.Range("E5").Value = "  +1.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.69"
$ws.Range("E6").Value = "  -1.12%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.141.90"
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  -0.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.157"
$ws.Range("E10").Value = "  -2.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.14"
$ws.Range("E11").Value = "  -1.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.498"
$ws.Range("E12").Value = "  -1.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000264"
$ws.Range("E13").Value = "  -1.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.99"
$ws.Range("E14").Value = "  -2.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.655.90"
$ws.Range("E15").Value = "  -0.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.756.89"
$ws.Range("E16").Value = "  -0.40%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.150.09"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.14"
$ws.Range("E18").Value = "  -0.95%  "
$ws.Range("E19").Value = "  -0.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "500.80"
$ws.Range("E20").Value = "  -1.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.32"
$ws.Range("E21").Value = "  +2.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.712"
$ws.Range("E22").Value = "  -3.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.01"
$ws.Range("E23").Value = "  -7.82%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.76"
$ws.Range("E24").Value = "  -0.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.48"
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.05"
$ws.Range("E27").Value = "  -0.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.91"
$ws.Range("E28").Value = "  +0.42%  "
$ws.Range("E29").Value = "  -0.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.81"
$ws.Range("E30").Value = "  +1.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "27.51"
$ws.Range("E31").Value = "  -1.34%  "
$ws.Range("E32").Value = "  -1.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.998"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.38"
$ws.Range("E34").Value = "  +1.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.45"
$ws.Range("E35").Value = "  -2.59%  "
$ws.Range("E36").Value = "  -1.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0893"
$ws.Range("E37").Value = "  +2.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "469.04"
$ws.Range("E38").Value = "  -0.97%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0418"
$ws.Range("E39").Value = "  -0.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.90"
$ws.Range("E40").Value = "  -4.12%  "
$ws.Range("E41").Value = "  +1.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.970.52"
$ws.Range("E42").Value = "  -4.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.116"
$ws.Range("E43").Value = "  -3.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.42"
$ws.Range("E44").Value = "  -2.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.280"
$ws.Range("E45").Value = "  -3.33%  "
$ws.Range("B46").Value = "PEPE"
$ws.Range("C46").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₃0604"
$ws.Range("E46").Value = "  +4.06%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.12"
$ws.Range("E47").Value = "  -3.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("E49").Value = "  -1.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.23"
$ws.Range("E50").Value = "  -4.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "119.44"
$ws.Range("E51").Value = "  -3.56%  "
